$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-8:
# the stored serial date value moves from 45184 to 45185 (one day later).
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value2 = 45185
}
